$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '59.021.07'
Set-TextCell 2 5 '  -4.61%  '
Set-TextCell 3 4 '2.517.60'
Set-TextCell 3 5 '  -2.27%  '
Set-TextCell 4 4 '0.998'
Set-TextCell 4 5 '  -0.23%  '
Set-TextCell 5 4 '535.94'
Set-TextCell 5 5 '  -2.45%  '
Set-TextCell 6 4 '144.95'
Set-TextCell 6 5 '  -6.18%  '
Set-TextCell 7 4 '0.996'
Set-TextCell 7 5 '  -0.34%  '
Set-TextCell 8 5 '  -3.13%  '
Set-TextCell 9 4 '2.517.19'
Set-TextCell 9 5 '  -2.55%  '
Set-TextCell 10 4 '0.0996'
Set-TextCell 10 5 '  -4.26%  '
Set-TextCell 11 5 '  -2.67%  '
Set-TextCell 12 4 '5.55'
Set-TextCell 12 5 '  -0.03%  '
Set-TextCell 13 5 '  -3.00%  '
Set-TextCell 14 4 '2.927.10'
Set-TextCell 14 5 '  -3.55%  '
Set-TextCell 15 4 '23.95'
Set-TextCell 15 5 '  -6.13%  '
Set-TextCell 16 4 '58.940.20'
Set-TextCell 16 5 '  -4.71%  '
Set-TextCell 17 5 '  -3.76%  '
Set-TextCell 18 4 '2.515.36'
Set-TextCell 18 5 '  -2.64%  '
Set-TextCell 19 4 '11.28'
Set-TextCell 19 5 '  -2.88%  '
Set-TextCell 20 5 '  -5.51%  '
Set-TextCell 21 4 '323.58'
Set-TextCell 21 5 '  -4.09%  '
Set-TextCell 22 4 '0.999'
Set-TextCell 22 5 '  +0.12%  '
Set-TextCell 23 4 '5.75'
Set-TextCell 23 5 '  -4.38%  '
Set-TextCell 24 4 '61.40'
Set-TextCell 24 5 '  -3.48%  '
Set-TextCell 25 5 '  -10.77%  '
Set-TextCell 26 2 'Kaspa'
Set-TextCell 26 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 26 4 '0.162'
Set-TextCell 26 5 '  -3.85%  '
Set-TextCell 27 2 'Binance-PegBSC-USD'
Set-TextCell 27 3 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 27 4 '0.996'
Set-TextCell 27 5 '  -0.26%  '
Set-TextCell 28 2 'WrappedeETH'
Set-TextCell 28 3 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell 28 4 '2.608.11'
Set-TextCell 28 5 '  -3.41%  '
Set-TextCell 29 4 '7.78'
Set-TextCell 29 5 '  -4.89%  '
Set-TextCell 30 4 '6.85'
Set-TextCell 30 5 '  -5.47%  '
Set-TextCell 31 4 '0.0₃0779'
Set-TextCell 31 5 '  -6.82%  '
Set-TextCell 32 4 '1.25'
Set-TextCell 32 5 '  -7.89%  '
Set-TextCell 33 5 '  -5.35%  '
Set-TextCell 34 4 '0.996'
Set-TextCell 34 5 '  -0.28%  '
Set-TextCell 35 4 '159.37'
Set-TextCell 35 5 '  -2.22%  '
Set-TextCell 36 4 '1.45'
Set-TextCell 36 5 '  +2.73%  '
Set-TextCell 37 4 '18.53'
Set-TextCell 37 5 '  -3.39%  '
Set-TextCell 38 4 '4.44'
Set-TextCell 38 5 '  -8.94%  '
Set-TextCell 39 5 '  -9.86%  '
Set-TextCell 40 4 '5.88'
Set-TextCell 40 5 '  -2.83%  '
Set-TextCell 41 4 '309.61'
Set-TextCell 41 5 '  -5.54%  '
Set-TextCell 42 4 '36.79'
Set-TextCell 42 5 '  -1.74%  '
Set-TextCell 43 4 '3.68'
Set-TextCell 43 5 '  -6.50%  '
Set-TextCell 44 4 '0.824'
Set-TextCell 44 5 '  -8.81%  '
Set-TextCell 45 4 '0.998'
Set-TextCell 45 5 '  -0.03%  '
Set-TextCell 46 4 '0.595'
Set-TextCell 46 5 '  -2.02%  '
Set-TextCell 47 5 '  -1.36%  '
Set-TextCell 48 4 '124.43'
Set-TextCell 48 5 '  +1.31%  '
Set-TextCell 49 4 '0.0931'
Set-TextCell 49 5 '  -3.61%  '
Set-TextCell 50 5 '  -4.65%  '
Set-TextCell 51 4 '0.0517'
Set-TextCell 51 5 '  -5.31%  '
